# Fill in the "Questionnaire" answers (Belbin self-perception inventory)
# so that each of the 7 sections sums to 10 points, matching the
# completed questionnaire committed upstream.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Questionnaire")

$ws.Range("B4").Value = 0
$ws.Range("B5").Value = 0
$ws.Range("B6").Value = 2
$ws.Range("B8").Value = 2
$ws.Range("B9").Value = 2
$ws.Range("B10").Value = 2
$ws.Range("B11").Value = 2

$ws.Range("B15").Value = 1
$ws.Range("B16").Value = 1
$ws.Range("B17").Value = 1
$ws.Range("B18").Value = 4
$ws.Range("B19").Value = 1
$ws.Range("B20").Value = 0
$ws.Range("B21").Value = 2
$ws.Range("B22").Value = 0

$ws.Range("B26").Value = 2
$ws.Range("B27").Value = 1
$ws.Range("B28").Value = 2
$ws.Range("B29").Value = 1
$ws.Range("B30").Value = 1
$ws.Range("B31").Value = 0
$ws.Range("B32").Value = 2
$ws.Range("B33").Value = 1

$ws.Range("B37").Value = 0
$ws.Range("B38").Value = 2
$ws.Range("B39").Value = 3
$ws.Range("B40").Value = 0
$ws.Range("B41").Value = 0
$ws.Range("B42").Value = 3
$ws.Range("B43").Value = 0
$ws.Range("B44").Value = 2

$ws.Range("B48").Value = 3
$ws.Range("B49").Value = 3
$ws.Range("B50").Value = 0
$ws.Range("B51").Value = 2
$ws.Range("B52").Value = 0
$ws.Range("B53").Value = 0
$ws.Range("B54").Value = 0
$ws.Range("B55").Value = 2

$ws.Range("B59").Value = 2
$ws.Range("B60").Value = 2
$ws.Range("B61").Value = 1
$ws.Range("B62").Value = 0
$ws.Range("B63").Value = 2
$ws.Range("B64").Value = 1
$ws.Range("B65").Value = 2
$ws.Range("B66").Value = 0

$ws.Range("B70").Value = 3
$ws.Range("B71").Value = 1
$ws.Range("B72").Value = 1
$ws.Range("B73").Value = 0
$ws.Range("B74").Value = 3
$ws.Range("B75").Value = 2
$ws.Range("B76").Value = 0
$ws.Range("B77").Value = 0

# Reflect the author's final on-screen state: cursor scrolled down on the
# Questionnaire sheet, and the "Grille d'évaluation" sheet becomes the
# active/selected tab.
$ws.Application.Goto($ws.Range("A49"))
$ws.Range("B76").Select()

$grille = $wb.Worksheets.Item("Grille d'évaluation")
$grille.Activate()
$grille.Range("P4:Q4").Select()
